$d = $word.ActiveDocument

# Phase 1: replace each original line with a unique placeholder token (by position)
# to avoid collisions, since several lines swap values with each other.
$rng = $d.Content
$null = $rng.Find.Execute("LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "01", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "02", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "03", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "04", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1006 -  Cálculo IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "05", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "06", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1037 -  Àlgebra Linear  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "07", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1053 -  Física III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "08", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3241 -  Química de Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "09", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1003 -  Cálculo I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "10", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1012 -  Estatística  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "11", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1018 -  Física I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "12", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1036 -  Geometria Analítica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "13", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1038 -  Física Experimental I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "14", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1039 -  Física Experimental III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "15", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1041 -  Física Experimental II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "16", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1042 -  Física Experimental IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "17", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1052 -  Cálculo III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "18", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "19", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3218 -  Introdução à Engenharia Física  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "20", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3260 -  Computação Científica em Python  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "21", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1004 -  Cálculo II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "22", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1019 -  Física II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "23", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOB1021 -  Física IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "24", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "25", 2)
$rng = $d.Content
$null = $rng.Find.Execute("LOM3236 -  Processos de Fabricação  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "26", 2)

# Phase 2: replace each placeholder with its final text value
$rng = $d.Content
$null = $rng.Find.Execute("01", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1053 -  Física III  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("02", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("03", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("04", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1036 -  Geometria Analítica  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("05", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1037 -  Àlgebra Linear  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("06", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1041 -  Física Experimental II  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("07", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1042 -  Física Experimental IV  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("08", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4095 -  Química Geral Experimental  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("09", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1039 -  Física Experimental III  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("10", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1018 -  Física I  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("11", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("12", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1004 -  Cálculo II  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("13", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1038 -  Física Experimental I  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("14", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1052 -  Cálculo III  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("15", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3236 -  Processos de Fabricação  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("16", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("17", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3218 -  Introdução à Engenharia Física  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("18", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1003 -  Cálculo I  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("19", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1006 -  Cálculo IV  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("20", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3241 -  Química de Materiais  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("21", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1021 -  Física IV  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("22", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("23", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3260 -  Computação Científica em Python  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("24", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("25", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1012 -  Estatística  (Requisito)", 2)
$rng = $d.Content
$null = $rng.Find.Execute("26", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1019 -  Física II  (Requisito)", 2)
